$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.658.76'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').Value = '3.408.34'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '568.14'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '157.04'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '3.412.74'
$ws.Range('E8').Value = '  -0.63%  '
$ws.Range('E9').Value = '  -7.47%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.24'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('E11').Value = '  -3.90%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.422'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -4.81%  '
$ws.Range('D13').Value = '3.997.55'
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '26.89'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -4.00%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000170'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -9.59%  '
$ws.Range('D17').Value = '63.729.48'
$ws.Range('E17').Value = '  -1.53%  '
$ws.Range('D18').Value = '3.374.03'
$ws.Range('E18').Value = '  -1.31%  '
$ws.Range('E19').Value = '  -4.69%  '
$ws.Range('E20').Value = '  -3.30%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '384.29'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.76'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.95%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '71.10'
$ws.Range('D24').Style = "Normal"
$ws.Range('E25').Value = '  -7.20%  '
$ws.Range('E26').Value = '  -5.57%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.67'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -5.42%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = "Normal"
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range('E31').Value = '  -7.14%  '
$ws.Range('E32').Value = '  -2.65%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '22.87'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.25%  '
$ws.Range('E35').Value = '  -4.53%  '
$ws.Range('E36').Value = '  -6.37%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '160.63'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.840'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +8.61%  '
$ws.Range('E39').Value = '  -4.73%  '
$ws.Range('D40').Value = '2.811.32'
$ws.Range('E40').Value = '  -3.17%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '25.89'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.53%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0723'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -5.22%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '42.99'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.46%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '6.38'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -9.43%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '25.60'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.81%  '
$ws.Range('E46').Value = '  -6.37%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0303'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.22%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.33'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +6.83%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '327.95'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.78%  '
$ws.Range('E50').Value = '  -5.16%  '
$ws.Range('E51').Value = '  -5.50%  '

Write-Host "Applied 85 cell changes"
